# "Generate Report for Archive"
#
# Updates the localization status report:
#   - Status "Ready for handoff" -> "In Translation" everywhere it appears
#     (Overview sheet columns E/F, and the Status column (C) on each
#     per-language sheet).
#   - The Status column(s) shrink to fit the new (shorter) text, so their
#     column width is reduced accordingly.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "In Translation"

# --- Overview sheet: zh-cn / de-de status columns (E & F), rows 2-4 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$usedOverview = $wsOverview.UsedRange
$lastRowOverview = $usedOverview.Rows.Count
for ($r = 2; $r -le $lastRowOverview; $r++) {
    if ($wsOverview.Cells.Item($r, 5).Value2 -eq $statusOld) {
        $wsOverview.Cells.Item($r, 5).Value = $statusNew
    }
    if ($wsOverview.Cells.Item($r, 6).Value2 -eq $statusOld) {
        $wsOverview.Cells.Item($r, 6).Value = $statusNew
    }
}

# --- Per-language sheets: Status column (C), rows 2-4 ---
$languageSheets = @("zh-cn", "de-de")
foreach ($sheetName in $languageSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count
    for ($r = 2; $r -le $lastRow; $r++) {
        if ($ws.Cells.Item($r, 3).Value2 -eq $statusOld) {
            $ws.Cells.Item($r, 3).Value = $statusNew
        }
    }
}

# --- Shrink the Status columns to fit the new, shorter text ---
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

foreach ($sheetName in $languageSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Columns.Item(3).ColumnWidth = 12.5
}
